$wb = $excel.ActiveWorkbook

# --- Update the "Conversión del día" note on sheet Hoja1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 1.84 = 6694.94 pesos`n✅ 6694.94 pesos = 1.83 = 947.86 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Update the rate figures on sheet tasas ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 544
$ws2.Range("O10").Value = 3642.05
$ws2.Range("N12").Value = 3651.7
$ws2.Range("O12").Value = 517.003
